# Apply "anglers data 1st take" edits to Table22.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric values in row 21 ---
$ws.Range("D21").Value = 802811
$ws.Range("E21").Value = 873327
$ws.Range("F21").Value = 728126
$ws.Range("G21").Value = 792618
$ws.Range("H21").Value = 880100
$ws.Range("I21").Value = 809100
$ws.Range("J21").Value = 748052
$ws.Range("K21").Value = 735674

# --- Update numeric values in row 22 ---
$ws.Range("B22").Value = 3880630
$ws.Range("C22").Value = 4294902
$ws.Range("E22").Value = 4442856
$ws.Range("F22").Value = 3778266
$ws.Range("G22").Value = 4075967
$ws.Range("H22").Value = 4382184
$ws.Range("I22").Value = 3875750
$ws.Range("J22").Value = 3539066
$ws.Range("K22").Value = 3539513

# --- Update the view state: scroll so column E is the leftmost visible column, and select K23 ---
$ws.Range("K23").Select()
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("K23").Select()
